# Append the "region" lookup table (rows 23-40) to Sheet1, mirroring the
# existing policy-hashtag table in columns I (hashtag) and K
# (CONCATENATE formula), paired up with a new column J of region names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Same 18 hashtags already used in I3:I20 / D3:D20, in the same order.
$hashtags = @(
    "#health_monitor",
    "#public_aw",
    "#other_pol",
    "#health_rsrc",
    "#task_fc",
    "#ext_border",
    "#int_border",
    "#mass_gath",
    "#gov_serv",
    "#emergency",
    "#schools",
    "#business",
    "#health_test",
    "#soc_dist",
    "#lockdown",
    "#curfew",
    "#hygeine",
    "#disinfo"
)

# Regions cycle in groups of 3 across the 18 rows.
$regions = @(
    "South Asia",
    "Europe & Central Asia",
    "Latin America & Caribbean"
)

$startRow = 23
for ($i = 0; $i -lt $hashtags.Length; $i++) {
    $row = $startRow + $i
    $hashtag = $hashtags[$i]
    $region = $regions[$i % 3]

    $ws.Cells.Item($row, 9).Value = $hashtag   # column I
    $ws.Cells.Item($row, 10).Value = $region   # column J
    $ws.Cells.Item($row, 11).Formula = "=CONCATENATE(""'"",I$row,""': '"",J$row,""',"")"  # column K
}

# Restore the view state captured in the saved workbook: scrolled so E13
# is the top-left visible cell, zoomed to 134%, with K23:K40 selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5    # column E
$win.ScrollRow = 13      # row 13
$win.Zoom = 134
[void]$ws.Range("K23:K40").Select()
